$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 737.5
$ws.Range("I18").Value = 866.6667
$ws.Range("J18").Value = 350
$ws.Range("K18").Value = 866.6667
$ws.Range("L18").Value = 350
$ws.Range("M18").Value = -582.6667
$ws.Range("N18").Value = -918

$ws.Range("H33").Value = 691.375
$ws.Range("I33").Value = 204.25
$ws.Range("J33").Value = 2152.75
$ws.Range("K33").Value = 204.25
$ws.Range("L33").Value = 2152.75
$ws.Range("M33").Value = 24.75

$ws.Range("H62").Value = 3969.4167
$ws.Range("I62").Value = 2878.9443
$ws.Range("J62").Value = 7240.8335
$ws.Range("K62").Value = 2878.9443
$ws.Range("L62").Value = 7240.8335
$ws.Range("M62").Value = -2254.9443

$ws.Range("H65").Value = 3969.4167
$ws.Range("I65").Value = 2878.9443
$ws.Range("J65").Value = 7240.8335
$ws.Range("K65").Value = 14394.7215
$ws.Range("L65").Value = 36204.1675
$ws.Range("M65").Value = -11274.7215

$ws.Range("H116").Value = 7059.45
$ws.Range("I116").Value = 6523.5
$ws.Range("J116").Value = 7416.75
$ws.Range("K116").Value = 6523.5
$ws.Range("L116").Value = 7416.75
$ws.Range("M116").Value = -3081.5
$ws.Range("N116").Value = -14300.75

$ws.Range("H121").Value = 2198.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2198.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6595.5
$ws.Range("N121").Value = -10089.5

$ws.Range("H125").Value = 988.4167
$ws.Range("I125").Value = 946.8333
$ws.Range("J125").Value = 1030
$ws.Range("K125").Value = 8521.4997
$ws.Range("L125").Value = 9270
$ws.Range("M125").Value = -6061.4997
$ws.Range("N125").Value = -14190

$ws.Range("H132").Value = 1054.091
$ws.Range("I132").Value = 976.0714
$ws.Range("J132").Value = 2692.5
$ws.Range("K132").Value = 2928.2142
$ws.Range("L132").Value = 8077.5
$ws.Range("M132").Value = -398.2142000000003

$ws.Range("H135").Value = 713.875
$ws.Range("I135").Value = 582.5238000000001
$ws.Range("J135").Value = 1633.3334
$ws.Range("K135").Value = 5242.7142
$ws.Range("L135").Value = 14700.0006
$ws.Range("M135").Value = -2707.7142
$ws.Range("N135").Value = -19770.0006

$ws.Range("H137").Value = 2499.577
$ws.Range("I137").Value = 662.875
$ws.Range("J137").Value = 3315.889
$ws.Range("K137").Value = 1988.625
$ws.Range("L137").Value = 9947.667000000001
$ws.Range("M137").Value = 561.375

$ws.Range("H138").Value = 2653.3696
$ws.Range("I138").Value = 2058.182
$ws.Range("J138").Value = 3198.9583
$ws.Range("K138").Value = 6174.545999999999
$ws.Range("L138").Value = 9596.874899999999
$ws.Range("M138").Value = -1034.545999999999
$ws.Range("N138").Value = -19876.8749


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15701.667
$ws.Range("I2").Value = 890
$ws.Range("J2").Value = 17553.125
$ws.Range("K2").Value = 890
$ws.Range("L2").Value = 17553.125
$ws.Range("M2").Value = -777
$ws.Range("N2").Value = -17779.125

$ws.Range("H40").Value = 23515
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 23515
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 23515
$ws.Range("N40").Value = -23867

$ws.Range("H107").Value = 45000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 45000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 45000
$ws.Range("N107").Value = -52680

$ws.Range("H110").Value = 2156.2632
$ws.Range("I110").Value = 1376.9286
$ws.Range("J110").Value = 4338.4
$ws.Range("K110").Value = 1376.9286
$ws.Range("L110").Value = 4338.4
$ws.Range("M110").Value = 668.0714

$ws.Range("H116").Value = 15701.667
$ws.Range("I116").Value = 890
$ws.Range("J116").Value = 17553.125
$ws.Range("K116").Value = 890
$ws.Range("L116").Value = 17553.125
$ws.Range("M116").Value = 1404
$ws.Range("N116").Value = -22141.125

$ws.Range("H124").Value = 41650
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 41650
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 41650
$ws.Range("N124").Value = -51470

$ws.Range("H132").Value = 3401.6316
$ws.Range("I132").Value = 2707.3572
$ws.Range("J132").Value = 5345.6
$ws.Range("K132").Value = 8122.071599999999
$ws.Range("L132").Value = 16036.8
$ws.Range("M132").Value = -5592.071599999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15701.667
$ws.Range("I3").Value = 890
$ws.Range("J3").Value = 17553.125
$ws.Range("K3").Value = 890
$ws.Range("L3").Value = 17553.125
$ws.Range("M3").Value = -776
$ws.Range("N3").Value = -17781.125

$ws.Range("H105").Value = 8156.222
$ws.Range("I105").Value = 3581.4
$ws.Range("J105").Value = 13874.75
$ws.Range("K105").Value = 3581.4
$ws.Range("L105").Value = 13874.75
$ws.Range("M105").Value = -1834.4
$ws.Range("N105").Value = -17368.75

$ws.Range("H107").Value = 973.6957
$ws.Range("I107").Value = 989.2632
$ws.Range("J107").Value = 899.75
$ws.Range("K107").Value = 989.2632
$ws.Range("L107").Value = 899.75
$ws.Range("M107").Value = 930.7368
$ws.Range("N107").Value = -4739.75

$ws.Range("H135").Value = 69999
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 69999
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1783.4
$ws.Range("I16").Value = 1640.8462
$ws.Range("J16").Value = 2710
$ws.Range("K16").Value = 1640.8462
$ws.Range("L16").Value = 2710
$ws.Range("M16").Value = -1353.8462
$ws.Range("N16").Value = -3284

$ws.Range("H31").Value = 29660.65
$ws.Range("I31").Value = 3494.4583
$ws.Range("J31").Value = 68909.94
$ws.Range("K31").Value = 3494.4583
$ws.Range("L31").Value = 68909.94
$ws.Range("M31").Value = -3199.4583
$ws.Range("N31").Value = -69499.94

$ws.Range("H34").Value = 29660.65
$ws.Range("I34").Value = 3494.4583
$ws.Range("J34").Value = 68909.94
$ws.Range("K34").Value = 3494.4583
$ws.Range("L34").Value = 68909.94
$ws.Range("M34").Value = -3292.4583
$ws.Range("N34").Value = -69313.94

$ws.Range("H105").Value = 7055.778
$ws.Range("I105").Value = 5748.5
$ws.Range("J105").Value = 9670.333000000001
$ws.Range("K105").Value = 5748.5
$ws.Range("L105").Value = 9670.333000000001
$ws.Range("M105").Value = -4001.5

$ws.Range("H107").Value = 1128.5294
$ws.Range("I107").Value = 1102.4
$ws.Range("J107").Value = 1165.8572
$ws.Range("K107").Value = 1102.4
$ws.Range("L107").Value = 1165.8572
$ws.Range("M107").Value = 817.5999999999999

$ws.Range("H111").Value = 91150.664
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 91150.664
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 91150.664
$ws.Range("N111").Value = -99330.664

$ws.Range("H113").Value = 1783.4
$ws.Range("I113").Value = 1640.8462
$ws.Range("J113").Value = 2710
$ws.Range("K113").Value = 1640.8462
$ws.Range("L113").Value = 2710
$ws.Range("M113").Value = 529.1538
$ws.Range("N113").Value = -7050

$ws.Range("H122").Value = 4818.3335
$ws.Range("I122").Value = 1438.0555
$ws.Range("J122").Value = 9888.75
$ws.Range("K122").Value = 4314.166499999999
$ws.Range("L122").Value = 29666.25
$ws.Range("M122").Value = -1864.166499999999
$ws.Range("N122").Value = -34566.25

$ws.Range("H125").Value = 59999
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 59999
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 59999
$ws.Range("N125").Value = -64919

$ws.Range("H132").Value = 6049.154
$ws.Range("I132").Value = 5888.4097
$ws.Range("J132").Value = 8500.5
$ws.Range("K132").Value = 17665.2291
$ws.Range("L132").Value = 25501.5
$ws.Range("M132").Value = -15135.2291


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6469.3184
$ws.Range("I5").Value = 729.4167
$ws.Range("J5").Value = 13357.2
$ws.Range("K5").Value = 2188.2501
$ws.Range("L5").Value = 40071.60000000001
$ws.Range("M5").Value = -2076.2501
$ws.Range("N5").Value = -40295.60000000001

$ws.Range("H131").Value = 14959155
$ws.Range("I131").Value = 2440.625
$ws.Range("J131").Value = 38889900
$ws.Range("K131").Value = 7321.875
$ws.Range("L131").Value = 116669700
$ws.Range("M131").Value = -2281.875
$ws.Range("N131").Value = -116679780

$ws.Range("H135").Value = 6469.3184
$ws.Range("I135").Value = 729.4167
$ws.Range("J135").Value = 13357.2
$ws.Range("K135").Value = 6564.7503
$ws.Range("L135").Value = 120214.8
$ws.Range("M135").Value = -4029.7503
$ws.Range("N135").Value = -125284.8


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13405.296
$ws.Range("I70").Value = 5450.0527
$ws.Range("J70").Value = 32299
$ws.Range("K70").Value = 5450.0527
$ws.Range("L70").Value = 32299
$ws.Range("M70").Value = -5180.0527

$ws.Range("H73").Value = 13405.296
$ws.Range("I73").Value = 5450.0527
$ws.Range("J73").Value = 32299
$ws.Range("K73").Value = 5450.0527
$ws.Range("L73").Value = 32299
$ws.Range("M73").Value = -4514.0527

$ws.Range("H113").Value = 6237.4
$ws.Range("I113").Value = 4461
$ws.Range("J113").Value = 7421.6665
$ws.Range("K113").Value = 4461
$ws.Range("L113").Value = 7421.6665
$ws.Range("M113").Value = -2291

$ws.Range("H120").Value = 44000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 44000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 44000
$ws.Range("N120").Value = -53676

$ws.Range("H124").Value = 74981.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 74981.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 74981.5
$ws.Range("N124").Value = -84801.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3198.9333
$ws.Range("I46").Value = 1998.3334
$ws.Range("J46").Value = 3499.0833
$ws.Range("K46").Value = 1998.3334
$ws.Range("L46").Value = 3499.0833
$ws.Range("M46").Value = -1810.3334


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74499.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 74499.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 74499.5
$ws.Range("N46").Value = -74961.5

$ws.Range("H51").Value = 23714.285
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 23714.285
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 23714.285
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -24734.285

$ws.Range("H125").Value = 82671.75
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 82671.75
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 82671.75
$ws.Range("N125").Value = -92511.75

$ws.Range("H134").Value = 74499.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 74499.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 223498.5
$ws.Range("N134").Value = -228568.5

